$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.465.86"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.645.20"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "300.51"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").Value = "0.3791"
$ws.Range("E7").Value = "  -1.15%  "
$ws.Range("D8").Value = "50.54"
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").Value = "0.3503"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D10").Value = "0.08061"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").Value = "1.217"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "22.11"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "6.293"
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("D15").Value = "7.254"
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("D16").Value = "0.00001212"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "1.642.29"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "95.31"
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("D19").Value = "0.06995"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "6.636"
$ws.Range("D21").Value = "17.41"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "12.44"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").Value = "23.461.49"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").Value = "2.415"
$ws.Range("E25").Value = "  -4.06%  "
$ws.Range("D26").Value = "2.993"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").Value = "21.03"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").Value = "151.69"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").Value = "131.67"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").Value = "1.828.20"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "6.865"
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("D33").Value = "2.141"
$ws.Range("E33").Value = "  -4.71%  "
$ws.Range("E34").Value = "  -7.24%  "
$ws.Range("D35").Value = "0.9885"
$ws.Range("E35").Value = "  -6.94%  "
$ws.Range("E36").Value = "  -3.95%  "
$ws.Range("D37").Value = "0.08784"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "5.913"
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("D39").Value = "0.2418"
$ws.Range("E39").Value = "  -3.21%  "
$ws.Range("D40").Value = "0.06786"
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("D41").Value = "12.88"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D42").Value = "0.6885"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").Value = "1.291"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("D44").Value = "15.49"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "0.6396"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "2.242"
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("D49").Value = "127.38"
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").Value = "0.07672"
$ws.Range("E50").Value = "  -2.97%  "
$ws.Range("D51").Value = "1.242"
$ws.Range("E51").Value = "  +3.09%  "
